$wb = $excel.ActiveWorkbook

# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on row 5 (the 7c3ff293-... handback entry) for both the zh-cn
# and de-de report sheets, reflecting the newly generated handback report.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-18 10:22:22"
$wsZhCn.Range("G5").Value = "2016-02-18 10:23:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-18 10:22:33"
$wsDeDe.Range("G5").Value = "2016-02-18 10:23:31"
